$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.725.54"
$ws.Range("E2").Value = "'  +0.51%  "
$ws.Range("D3").Value = "'2.559.88"
$ws.Range("E3").Value = "'  +0.51%  "
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'302.65"
$ws.Range("E5").Value = "'  +2.52%  "
$ws.Range("D6").Value = "'97.14"
$ws.Range("E6").Value = "'  +7.26%  "
$ws.Range("E7").Value = "'  +0.72%  "
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("E9").Value = "'  +0.40%  "
$ws.Range("D10").Value = "'36.56"
$ws.Range("E10").Value = "'  +3.60%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "'  +1.12%  "
$ws.Range("E12").Value = "'  +9.67%  "
$ws.Range("D13").Value = "'7.68"
$ws.Range("E13").Value = "'  +1.44%  "
$ws.Range("D14").Value = "'2.580.21"
$ws.Range("E14").Value = "'  +1.58%  "
$ws.Range("E15").Value = "'  +2.58%  "
$ws.Range("D16").Value = "'14.57"
$ws.Range("E16").Value = "'  +4.01%  "
$ws.Range("D17").Value = "'42.797.64"
$ws.Range("E17").Value = "'  +0.60%  "
$ws.Range("D18").Value = "'13.63"
$ws.Range("E18").Value = "'  +9.00%  "
$ws.Range("D19").Value = "'0.0₃0988"
$ws.Range("E19").Value = "'  +2.83%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "'  +0.24%  "
$ws.Range("D21").Value = "'71.61"
$ws.Range("E21").Value = "'  -0.84%  "
$ws.Range("D22").Value = "'256.71"
$ws.Range("E22").Value = "'  -0.15%  "
$ws.Range("D23").Value = "'2.94"
$ws.Range("E23").Value = "'  +2.81%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "'  -0.18%  "
$ws.Range("D25").Value = "'28.17"
$ws.Range("E25").Value = "'  -3.94%  "
$ws.Range("D27").Value = "'39.27"
$ws.Range("E27").Value = "'  +9.91%  "
$ws.Range("D28").Value = "'10.09"
$ws.Range("E28").Value = "'  +1.84%  "
$ws.Range("E29").Value = "'  -0.13%  "
$ws.Range("E30").Value = "'  +2.35%  "
$ws.Range("D31").Value = "'156.21"
$ws.Range("E31").Value = "'  +3.94%  "
$ws.Range("E32").Value = "'  +1.15%  "
$ws.Range("E33").Value = "'  +1.26%  "
$ws.Range("D34").Value = "'27.26"
$ws.Range("E34").Value = "'  +12.79%  "
$ws.Range("D35").Value = "'3.35"
$ws.Range("E35").Value = "'  -0.21%  "
$ws.Range("E36").Value = "'  +2.20%  "
$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = "'  +2.61%  "
$ws.Range("D38").Value = "'18.16"
$ws.Range("E38").Value = "'  +16.99%  "
$ws.Range("E39").Value = "'  +0.86%  "
$ws.Range("D40").Value = "'3.85"
$ws.Range("E40").Value = "'  +2.28%  "
$ws.Range("B41").Value = "'NEARProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'3.37"
$ws.Range("E41").Value = "'  -0.48%  "
$ws.Range("B42").Value = "'ApeXProtocol"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'2.04"
$ws.Range("E42").Value = "'  +28.66%  "
$ws.Range("D43").Value = "'0.0306"
$ws.Range("E43").Value = "'  -0.20%  "
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "'  +0.04%  "
$ws.Range("B45").Value = "'Maker"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.061.02"
$ws.Range("E45").Value = "'  +0.08%  "
$ws.Range("D46").Value = "'88.48"
$ws.Range("E46").Value = "'  +5.27%  "
$ws.Range("D47").Value = "'9.30"
$ws.Range("E47").Value = "'  +6.81%  "
$ws.Range("D48").Value = "'77.23"
$ws.Range("E48").Value = "'  +12.61%  "
$ws.Range("D49").Value = "'2.810.64"
$ws.Range("E49").Value = "'  +0.79%  "
$ws.Range("D50").Value = "'103.99"
$ws.Range("E50").Value = "'  +1.13%  "
$ws.Range("E51").Value = "'  +3.33%  "
